$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.062521815299988
$ws.Range("B1").Value = 3.679196834564209
$ws.Range("C1").Value = 3.297560453414917
$ws.Range("D1").Value = 2.003660440444946
$ws.Range("E1").Value = 1.155933141708374
